$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 82
$ws.Cells.Item(82, 8).Value = 2026.3334
$ws.Cells.Item(82, 9).Value = 2026.3334
$ws.Cells.Item(82, 11).Value = 6079.0002
$ws.Cells.Item(82, 13).Value = -5673.0002
# Row 85
$ws.Cells.Item(85, 8).Value = 2026.3334
$ws.Cells.Item(85, 9).Value = 2026.3334
$ws.Cells.Item(85, 11).Value = 6079.0002
$ws.Cells.Item(85, 13).Value = -4675.0002
# Row 113
$ws.Cells.Item(113, 8).Value = 90913080
$ws.Cells.Item(113, 9).Value = 142859420
$ws.Cells.Item(113, 10).Value = 6972
$ws.Cells.Item(113, 11).Value = 142859420
$ws.Cells.Item(113, 12).Value = 6972
$ws.Cells.Item(113, 13).Value = -142856166
$ws.Cells.Item(113, 14).Value = -13480
# Row 125
$ws.Cells.Item(125, 8).Value = 500
$ws.Cells.Item(125, 9).Value = 500
$ws.Cells.Item(125, 11).Value = 4500
$ws.Cells.Item(125, 13).Value = -2040
# Row 129
$ws.Cells.Item(129, 8).Value = 889.6429000000001
$ws.Cells.Item(129, 9).Value = 616.75
$ws.Cells.Item(129, 10).Value = 910.63464
$ws.Cells.Item(129, 11).Value = 1850.25
$ws.Cells.Item(129, 12).Value = 2731.90392
$ws.Cells.Item(129, 13).Value = 3149.75
$ws.Cells.Item(129, 14).Value = -12731.90392
# Row 138
$ws.Cells.Item(138, 8).Value = 2848.6758
$ws.Cells.Item(138, 9).Value = 1871.8
$ws.Cells.Item(138, 10).Value = 3514.7273
$ws.Cells.Item(138, 11).Value = 5615.4
$ws.Cells.Item(138, 12).Value = 10544.1819
$ws.Cells.Item(138, 13).Value = -475.3999999999996
$ws.Cells.Item(138, 14).Value = -20824.1819

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 35
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 13).ClearContents()
# Row 61
$ws.Cells.Item(61, 8).Value = 4672.4546
$ws.Cells.Item(61, 9).Value = 4499.7334
$ws.Cells.Item(61, 10).Value = 5042.5713
$ws.Cells.Item(61, 11).Value = 4499.7334
$ws.Cells.Item(61, 12).Value = 5042.5713
$ws.Cells.Item(61, 13).Value = -4287.7334
$ws.Cells.Item(61, 14).Value = -5466.5713
# Row 136
$ws.Cells.Item(136, 8).Value = 4672.4546
$ws.Cells.Item(136, 9).Value = 4499.7334
$ws.Cells.Item(136, 10).Value = 5042.5713
$ws.Cells.Item(136, 11).Value = 13499.2002
$ws.Cells.Item(136, 12).Value = 15127.7139
$ws.Cells.Item(136, 13).Value = -10949.2002
$ws.Cells.Item(136, 14).Value = -20227.7139

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Cells.Item(36, 8).Value = 2818.5
$ws.Cells.Item(36, 9).Value = 2818.5
$ws.Cells.Item(36, 11).Value = 2818.5
$ws.Cells.Item(36, 13).Value = -2284.5
# Row 105
$ws.Cells.Item(105, 8).Value = 2501668.8
$ws.Cells.Item(105, 9).Value = 1597.8182
$ws.Cells.Item(105, 11).Value = 1597.8182
$ws.Cells.Item(105, 13).Value = 149.1818000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 88.42856999999999
$ws.Cells.Item(7, 9).Value = 82.25
$ws.Cells.Item(7, 10).Value = 96.666664
$ws.Cells.Item(7, 11).Value = 82.25
$ws.Cells.Item(7, 12).Value = 96.666664
$ws.Cells.Item(7, 13).Value = 30.75
$ws.Cells.Item(7, 14).Value = -322.666664
# Row 16
$ws.Cells.Item(16, 8).Value = 1212
$ws.Cells.Item(16, 9).Value = 1101.1428
$ws.Cells.Item(16, 10).Value = 1600
$ws.Cells.Item(16, 11).Value = 1101.1428
$ws.Cells.Item(16, 12).Value = 1600
$ws.Cells.Item(16, 13).Value = -814.1428000000001
$ws.Cells.Item(16, 14).Value = -2174
# Row 31
$ws.Cells.Item(31, 8).Value = 1896.1915
$ws.Cells.Item(31, 9).Value = 874.7353000000001
$ws.Cells.Item(31, 10).Value = 4567.6924
$ws.Cells.Item(31, 11).Value = 874.7353000000001
$ws.Cells.Item(31, 12).Value = 4567.6924
$ws.Cells.Item(31, 13).Value = -579.7353000000001
$ws.Cells.Item(31, 14).Value = -5157.6924
# Row 34
$ws.Cells.Item(34, 8).Value = 1896.1915
$ws.Cells.Item(34, 9).Value = 874.7353000000001
$ws.Cells.Item(34, 10).Value = 4567.6924
$ws.Cells.Item(34, 11).Value = 874.7353000000001
$ws.Cells.Item(34, 12).Value = 4567.6924
$ws.Cells.Item(34, 13).Value = -672.7353000000001
$ws.Cells.Item(34, 14).Value = -4971.6924
# Row 113
$ws.Cells.Item(113, 8).Value = 1212
$ws.Cells.Item(113, 9).Value = 1101.1428
$ws.Cells.Item(113, 10).Value = 1600
$ws.Cells.Item(113, 11).Value = 1101.1428
$ws.Cells.Item(113, 12).Value = 1600
$ws.Cells.Item(113, 13).Value = 1068.8572
$ws.Cells.Item(113, 14).Value = -5940
# Row 134
$ws.Cells.Item(134, 8).Value = 802.1818
$ws.Cells.Item(134, 9).Value = 759
$ws.Cells.Item(134, 10).Value = 1234
$ws.Cells.Item(134, 11).Value = 2277
$ws.Cells.Item(134, 12).Value = 3702
$ws.Cells.Item(134, 13).Value = 258
$ws.Cells.Item(134, 14).Value = -8772

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Cells.Item(64, 8).Value = 2459.6667
$ws.Cells.Item(64, 9).Value = 1666.6666
$ws.Cells.Item(64, 11).Value = 4999.9998
$ws.Cells.Item(64, 13).Value = -4729.9998
# Row 67
$ws.Cells.Item(67, 8).Value = 2459.6667
$ws.Cells.Item(67, 9).Value = 1666.6666
$ws.Cells.Item(67, 11).Value = 4999.9998
$ws.Cells.Item(67, 13).Value = -4063.9998
# Row 122
$ws.Cells.Item(122, 8).Value = 882.4
$ws.Cells.Item(122, 9).Value = 454
$ws.Cells.Item(122, 10).Value = 1168
$ws.Cells.Item(122, 11).Value = 4086
$ws.Cells.Item(122, 12).Value = 10512
$ws.Cells.Item(122, 13).Value = -1636
$ws.Cells.Item(122, 14).Value = -15412
# Row 131
$ws.Cells.Item(131, 8).Value = 799.98
$ws.Cells.Item(131, 10).Value = 825.9895
$ws.Cells.Item(131, 12).Value = 2477.9685
$ws.Cells.Item(131, 14).Value = -12557.9685
# Row 132
$ws.Cells.Item(132, 8).Value = 1491.5
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 1491.5
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 13423.5
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -18483.5
# Row 137
$ws.Cells.Item(137, 8).Value = 8831.177
$ws.Cells.Item(137, 9).Value = 50249.5
$ws.Cells.Item(137, 10).Value = 3308.7334
$ws.Cells.Item(137, 11).Value = 150748.5
$ws.Cells.Item(137, 12).Value = 9926.200199999999
$ws.Cells.Item(137, 13).Value = -145648.5
$ws.Cells.Item(137, 14).Value = -20126.2002

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 3968.1
$ws.Cells.Item(7, 9).Value = 4090.3125
$ws.Cells.Item(7, 10).Value = 3479.25
$ws.Cells.Item(7, 11).Value = 4090.3125
$ws.Cells.Item(7, 12).Value = 3479.25
$ws.Cells.Item(7, 13).Value = -3978.3125
$ws.Cells.Item(7, 14).Value = -3703.25
# Row 126
$ws.Cells.Item(126, 8).Value = 3968.1
$ws.Cells.Item(126, 9).Value = 4090.3125
$ws.Cells.Item(126, 10).Value = 3479.25
$ws.Cells.Item(126, 11).Value = 12270.9375
$ws.Cells.Item(126, 12).Value = 10437.75
$ws.Cells.Item(126, 13).Value = -9800.9375
$ws.Cells.Item(126, 14).Value = -15377.75
# Row 139
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).ClearContents()
# Row 89
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).ClearContents()
# Row 139
$ws.Cells.Item(139, 8).Value = 52714.668
$ws.Cells.Item(139, 9).Value = 52714
$ws.Cells.Item(139, 11).Value = 52714
$ws.Cells.Item(139, 13).Value = -47574

Write-Host "Edits applied successfully"
